$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new negative test-case row (row 4) for the Header module, mirroring
# row 2's formatting but with a distinct candidate name ("Prasannaee") and an
# invalid/short salary value (900) to exercise the negative-test path.
$ws.Range("C2:P2").Copy()
$ws.Range("C4:P4").PasteSpecial(-4122)
$ws.Range("K4").ClearFormats()

$ws.Range("C4").Value = "Teaching"
$ws.Range("D4").Value = "Teacher"
$ws.Range("E4").Value = "Prasannaee"
$ws.Range("F4").Value = "siva"
$ws.Range("G4").Value = 900
$ws.Range("H4").Value = 35822
$ws.Range("I4").Value = "Male"
$ws.Range("J4").Value = 9876543210
$ws.Range("K4").Value = "prasanna565@gmail.com"
$ws.Range("L4").Value = "Teacher"
$ws.Range("M4").Value = 45563
$ws.Range("N4").Value = "Tamil Nadu"
$ws.Range("O4").Value = "Arcot"
$ws.Range("P4").Value = 654329

# Move the active selection onto the newly added row, as Excel would after entry
$ws.Range("C4:P4").Select()
